$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 22 in the "Logs" sheet: Testmail #7 - a return-request mail that was
# auto-answered (Retour / Terugbetaling category).
$ws.Cells.Item(22, 1).Value = "Hoe kan ik iets retourneren?"
$ws.Cells.Item(22, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(22, 3).Value = "Testmail #7: Hoe kan ik iets retourneren?"
$ws.Cells.Item(22, 4).Value = "Retour / Terugbetaling"
$ws.Cells.Item(22, 5).Value = "Beste klant,`nBedankt voor je bericht. Als je een artikel wilt retourneren, kun je dit doen door contact op te nemen met onze klantenservice via support@bedrijf.nl. Zij zullen je verder begeleiden bij het retourproces en eventuele vragen beantwoorden.`nBedankt voor je begrip en medewerking.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$ws.Cells.Item(22, 6).Value = "2025-06-29 14:55:39"
$ws.Cells.Item(22, 7).Value = "Ja"
$ws.Cells.Item(22, 8).Value = "Nee"
$ws.Cells.Item(22, 9).Value = "Ja"

# Grow the conditional-formatting ranges (D/G/H/I) one row down so they keep
# covering the whole data range now that row 22 exists.
$ws.Range("D2:D21").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D22"))
$ws.Range("G2:G21").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G22"))
$ws.Range("H2:H21").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H22"))
$ws.Range("I2:I21").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I22"))

# Dashboard summary: "Retour / Terugbetaling" tally goes from 2 to 3.
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Range("B5").Value = 3
